# Row column count and check for any fields in row
#
# - Adds a new "RowColumnTable" worksheet at the end with an expected-field
#   label and a sample data row ("Sonya Frost ...").
# - Updates the "Select Input" sheet's "All selected colors are : " message
#   to include the actually-selected color ("Green").
# - Updates the remembered selection/active-cell for a few sheets and makes
#   "Select Input" the active (visible) tab instead of "Simple Form Demo".

$wb = $excel.ActiveWorkbook

$wsSelectInput  = $wb.Worksheets.Item("Select Input")
$wsRadioButtons = $wb.Worksheets.Item("Radio Buttons Demo")
$wsSimpleForm   = $wb.Worksheets.Item("Simple Form Demo")

# --- Add the new "RowColumnTable" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRowColumnTable = $wb.Worksheets.Add($null, $lastSheet)
$wsRowColumnTable.Name = "RowColumnTable"

# Populate it: row 2 (the data row) is written before row 1 (the header
# label) so the shared-string table ends up with the sample data string
# immediately after the existing strings, followed by the header label.
$wsRowColumnTable.Range("A2").Value = "Sonya Frost Software Engineer Edinburgh 23 2008/12/13 `$103,600"
$wsRowColumnTable.Range("A1").Value = "ExpectedField"

# --- Update the "All selected colors are : " message on Select Input ---
$wsSelectInput.Range("B3").Value = "All selected colors are : Green"

# --- Update remembered selections on each sheet ---
[void]$wsRadioButtons.Activate()
[void]$wsRadioButtons.Range("A4").Select()

[void]$wsSimpleForm.Activate()
[void]$wsSimpleForm.Range("C1").Select()

[void]$wsRowColumnTable.Activate()
[void]$wsRowColumnTable.Range("B6").Select()

# "Select Input" ends up as the active/visible tab with B3 selected.
[void]$wsSelectInput.Activate()
[void]$wsSelectInput.Range("B3").Select()
